$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the reported p-value for row 8 (column J, p_reported) to reflect
# the replicated result.
$ws.Cells.Item(8, 10).Value = "<0.02"

# Insert a new column before column K (11th column) to hold the new
# "p_replicated" field, shifting the existing N/published/ci_lower/ci_upper
# columns one to the right.
$ws.Columns.Item(11).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 11).Value = "p_replicated"

# Match the best-fit column width Excel would compute for the new header
# text ("p_replicated"), stored in the XML as width 11.5.
$ws.Columns.Item(11).ColumnWidth = 10.666666666666666
